# Refresh the cryptocurrency "Price" (D) and "Volume(1h)" (E) columns with
# freshly scraped figures, matching the GitHub Actions daily-update commit.
#
# Price text is written with a leading apostrophe (quote-prefix) whenever it
# would otherwise be auto-parsed as a number (e.g. "0.999", "316.34"), so
# that Excel keeps it as literal text - exactly like the original cells
# (e.g. "1.00", "95.07"). The Style reset back to "Normal" right afterwards
# clears the quote-prefix cell format flag Excel applies automatically,
# so the cell keeps its original (default) formatting - only its text value
# changes. Values that can't be parsed as numbers anyway (e.g. "42.834.57",
# which contains two periods) are written plainly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.834.57"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "2.527.08"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'316.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "'94.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("D7").Value = "'0.579"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").Value = "'0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.43%  "
$ws.Range("D10").Value = "'35.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'0.0811"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "'7.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").Value = "'0.110"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "2.916.28"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "2.523.42"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'15.23"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "'0.847"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.71%  "
$ws.Range("D18").Value = "42.915.77"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'12.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'6.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.49%  "
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "'69.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").Value = "'251.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").Value = "'2.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("D26").Value = "'26.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +3.43%  "
$ws.Range("D29").Value = "'39.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.47%  "
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'154.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  +2.60%  "
$ws.Range("D34").Value = "'19.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "'0.0790"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'2.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.81%  "
$ws.Range("E38").Value = "  -2.62%  "
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("D40").Value = "'23.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$ws.Range("D41").Value = "'2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.26%  "
$ws.Range("D42").Value = "'0.0304"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "'3.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.95%  "
$ws.Range("D46").Value = "2.016.68"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").Value = "'85.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("D48").Value = "'8.79"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "2.770.72"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").Value = "'73.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("D51").Value = "'102.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
